{"js": "// Word JS API (Office.js) edit script\n// This reproduces the diff: several filled-in template values are cleared\n// out (name, NIK, birthplace/date, nationality/religion, job, address,\n// scholarship purpose) because the doc now renders via Pug instead of EJS,\n// and the \"Tommo, 29 Desember 2021\" signature date is bumped to 30.\nconst body = context.document.body;\n\nasync function replaceOnce(find, replacement) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + find);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nasync function replaceNth(find, index, replacement) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length <= index) {\n    throw new Error(\"No match #\" + index + \" found for: \" + find);\n  }\n  results.items[index].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// \"Nama\" ... \": Putu Mahendra\" -> \": \"\nawait replaceOnce(\": Putu Mahendra\", \": \");\n\n// \"NIK\" ... \" 7602111711990001\" -> \" \"\nawait replaceOnce(\" 7602111711990001\", \" \");\n\n// \"Tempat, Tanggal Lahir\" ... \": Tommo\" (1st occurrence) -> \": \"\nawait replaceNth(\": Tommo\", 0, \": \");\n\n// same paragraph ... \", 2021-12-09\" -> \", \"\nawait replaceOnce(\", 2021-12-09\", \", \");\n\n// \"Warganegara/Agama\" ... \": WNI/Hindu\" -> \": /\"\nawait replaceOnce(\": WNI/Hindu\", \": /\");\n\n// \"Pekerjaan\" ... \": Mahasiswa\" -> \": \"\nawait replaceOnce(\": Mahasiswa\", \": \");\n\n// \"Alamat\" ... \": Tommo\" (2nd/last occurrence) -> \": undefined\"\nawait replaceNth(\": Tommo\", 0, \": undefined\");\n\n// \"... untuk Beasiswa\" -> \"... untuk \" (drop \" Beasiswa\")\nawait replaceOnce(\" Beasiswa\", \" \");\n\n// Signature date bump\nawait replaceOnce(\"Tommo, 29 Desember 2021\", \"Tommo, 30 Desember 2021\");\n", "ps1": "# Word COM interop edit script\n# Reproduces the diff: several filled-in template values are cleared out\n# (name, NIK, birthplace/date, nationality/religion, job, address,\n# scholarship purpose) because the doc now renders via Pug instead of EJS,\n# and the \"Tommo, 29 Desember 2021\" signature date is bumped to 30.\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\nfunction Replace-InRange($range, [string]$findText, [string]$replaceText) {\n    $ok = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $ok) {\n        throw \"Find failed for: $findText\"\n    }\n}\n\n# Paragraph 13 - \"Nama\\t\\t: Putu Mahendra\" -> \"Nama\\t\\t: \"\nReplace-InRange $paras.Item(13).Range \": Putu Mahendra\" \": \"\n\n# Paragraph 14 - \"NIK\\t\\t: 7602111711990001\" -> \"NIK\\t\\t: \"\nReplace-InRange $paras.Item(14).Range \" 7602111711990001\" \" \"\n\n# Paragraph 16 - \"Tempat, Tanggal Lahir\\t: Tommo, 2021-12-09\" -> \"...\\t: , \"\nReplace-InRange $paras.Item(16).Range \": Tommo\" \": \"\nReplace-InRange $paras.Item(16).Range \", 2021-12-09\" \", \"\n\n# Paragraph 17 - \"Warganegara/Agama\\t\\t: WNI/Hindu\" -> \"...\\t\\t: /\"\nReplace-InRange $paras.Item(17).Range \": WNI/Hindu\" \": /\"\n\n# Paragraph 18 - \"Pekerjaan\\t\\t: Mahasiswa\" -> \"...\\t\\t: \"\nReplace-InRange $paras.Item(18).Range \": Mahasiswa\" \": \"\n\n# Paragraph 20 - \"Alamat\\t\\t: Tommo\" -> \"...\\t\\t: undefined\"\nReplace-InRange $paras.Item(20).Range \": Tommo\" \": undefined\"\n\n# Paragraph 23 - \"Surat Keterangan ini dibuat untuk Beasiswa.\" -> \"...untuk .\"\nReplace-InRange $paras.Item(23).Range \" Beasiswa\" \" \"\n\n# Paragraph 28 - \"Mengetahui,\\tTommo, 29 Desember 2021\" -> \"...Tommo, 30 Desember 2021\"\nReplace-InRange $paras.Item(28).Range \"Tommo, 29 Desember 2021\" \"Tommo, 30 Desember 2021\"\n"}
